# RPA datasets push 2024-04-23
# Refresh the underwriting dataset on Sheet1:
#   - remove the BNK / 비엔케이제2호스팩 row
#   - remove the 유진 / 유진스팩10호 row
#   - insert a new 신한제13호스팩 row just before the 신한제12호스팩 row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Delete the BNK / 비엔케이제2호스팩 row (data row right under the header).
$ws.Rows.Item(2).EntireRow.Delete()

# 2) Delete the 유진 / 유진스팩10호 row. After step 1, the table shifted up by
#    one row, so this row is now row 10 (was row 11 originally).
$ws.Rows.Item(10).EntireRow.Delete()

# 3) Insert a brand-new row for 신한제13호스팩 right above 신한제12호스팩
#    (now row 8 after the two deletions above) and populate it.
$ws.Rows.Item(8).EntireRow.Insert()

$newRow = 8
$ws.Cells.Item($newRow, 1).Value = "신한"

# Columns B, F and G hold dates formatted as plain text (e.g. "2024-04-11"),
# just like every other row in the table. Force text formatting before
# assigning so Excel doesn't auto-convert the literal into a date serial,
# then drop back to the workbook's default (unstyled) cell style so the
# new cells look just like their neighbours.
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).Value = "2024-04-11"
$ws.Cells.Item($newRow, 2).Style = "Normal"

$ws.Cells.Item($newRow, 3).Value = "신한제13호스팩"
$ws.Cells.Item($newRow, 4).Value = "신한"
$ws.Cells.Item($newRow, 5).Value = "신한"

$ws.Cells.Item($newRow, 6).NumberFormat = "@"
$ws.Cells.Item($newRow, 6).Value = "2024-04-15"
$ws.Cells.Item($newRow, 6).Style = "Normal"

$ws.Cells.Item($newRow, 7).NumberFormat = "@"
$ws.Cells.Item($newRow, 7).Value = "2024-04-22"
$ws.Cells.Item($newRow, 7).Style = "Normal"

$ws.Cells.Item($newRow, 8).Value = 6000
$ws.Cells.Item($newRow, 9).Value = 3000000
$ws.Cells.Item($newRow, 10).Value = 2000
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 100
